$d = $word.ActiveDocument

function Set-ParagraphRuns($paraIndex, $innerXml) {
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    $sub = $d.Range($r.Start, $r.End - 1)
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $innerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $sub.InsertXML($pkg)
}

# Paragraph 2: "Install " + "DB-Configurator" + " using SAG-installer." -> single run
Set-ParagraphRuns 2 '<w:r><w:t>Install DB-Configurator using SAG-installer.</w:t></w:r>'

# Paragraph 3: merge first 4 runs ("After installation ... under " + "folder" + " " + "~") into one,
# keep "\common\db" and " : bin, conf, lib, logs" as separate trailing runs.
Set-ParagraphRuns 3 '<w:r><w:t>After installation many new folders are visible under folder ~</w:t></w:r><w:r><w:t>\common\db</w:t></w:r><w:r><w:t xml:space="preserve"> : bin, conf, lib, logs</w:t></w:r>'

# Paragraph 4: merge first 4 runs ("Execute" + " dbconfigurator" + ": " + " double click on file ") into one,
# keep the bold "dbConfiguratorUI" run and the trailing runs unchanged.
Set-ParagraphRuns 4 '<w:r><w:t xml:space="preserve">Execute dbconfigurator:  double click on file </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>dbConfiguratorUI</w:t></w:r><w:r><w:t xml:space="preserve"> under ~</w:t></w:r><w:r><w:t>\common\db\bin</w:t></w:r><w:r><w:t>.</w:t></w:r>'

# Paragraph 7: "Create DataBase(Is_DB) in MySql" + "." -> single run
Set-ParagraphRuns 7 '<w:r><w:t>Create DataBase(Is_DB) in MySql.</w:t></w:r>'

# Paragraph 26: "If you get " + "Database Connection error: ..." -> single run
Set-ParagraphRuns 26 '<w:r><w:t xml:space="preserve">If you get Database Connection error: Access denied for user ''root''@''localhost'' (using password: YES)</w:t></w:r>'
